$d = $word.ActiveDocument

# =====================================================================
# Step 1: Heading text change.
#   "Use case "Watch detailed infomation of available bikes" "
#    -> "Use case "Select Bike" "
#   authored as three separate (identically formatted) runs.
# =====================================================================

$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*detailed infomation of available bikes*") {
        $headingPara = $p
        break
    }
}

$hr = $headingPara.Range
$start = $hr.Start

$find = $hr.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$oldHeading = "Use case " + [char]0x201C + "Watch detailed infomation of available bikes" + [char]0x201D + " "
$newHeading = "Use case " + [char]0x201C + "Select Bike" + [char]0x201D + " "
$find.Execute($oldHeading, $false, $false, $false, $false, $false, $true, 1, $false, $newHeading, 2) | Out-Null

# Force the merged text back into three distinct runs (all share the same
# rPr, so we nudge Bold off/on across each inner boundary to stop the
# engine from re-coalescing them).
$seg1Len = ("Use case " + [char]0x201C).Length
$seg2Len = "Select Bike".Length

$r1 = $d.Range($start, $start + $seg1Len)
$r1.Font.Bold = $false
$r1.Font.Bold = $true

$r2 = $d.Range($start + $seg1Len, $start + $seg1Len + $seg2Len)
$r2.Font.Bold = $false
$r2.Font.Bold = $true

# =====================================================================
# Step 2: Move the _GoBack bookmark.
#   Before: sits at the end of the "Buoc 2: ... ve xe." paragraph.
#   After:  sits at the end of the "... xem thong tin." paragraph
#           (i.e. the paragraph right before it).
# =====================================================================

$paraTin = $null
$paraXe = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*xem th*ng tin.*") {
        $paraTin = $p
    }
    if ($p.Range.Text -like "*ra c*c th*ng tin v*xe.*") {
        $paraXe = $p
    }
}

# The target insertion point is the end of $paraTin's content, i.e. right
# before its paragraph mark. Adding a bookmark with Bookmarks.Add on a
# zero-width range that sits exactly on a paragraph-end boundary is
# mis-handled, so first nudge the boundary away by inserting a temporary
# marker character, add the bookmark next to it, then remove the marker.
$endTin = $paraTin.Range.End - 1
$tmp = $d.Range($endTin, $endTin)
$tmp.InsertAfter("@")

$bmRange = $d.Range($endTin, $endTin)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$d.Range($endTin, $endTin + 1).Delete()
